$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 5) with a new user record, mirroring the
# existing firstName/lastName/email/password columns.
$ws.Range("A5").Value = "baryalai"
$ws.Range("B5").Value = "joyan"
$ws.Range("C5").Value = "user@test.com"
$ws.Range("D5").Value = "user2"

# C5 gets the same mailto: hyperlink treatment as the other email cells
# above it (C2:C4).
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:user@test.com")
$ws.Range("C5").Style = $ws.Range("C4").Style

# Move the active selection to A3 (matches the saved selection state).
$ws.Range("A3").Select() | Out-Null
